$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# Normalize particle casing in state/municipality names (and one all-caps fix)
$ws.Range('B6').Value = 'Rincón De Romos'
$ws.Range('B19').Value = 'Comitán De Domínguez'
$ws.Range('B26').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B28').Value = 'San Cristóbal De Las Casas'
$ws.Range('B64').Value = 'Guadalupe Y Calvo'
$ws.Range('B67').Value = 'Hidalgo Del Parral'
$ws.Range('B93').Value = 'San Francisco De Borja'
$ws.Range('B94').Value = 'San Francisco De Conchos'
$ws.Range('B95').Value = 'San Francisco Del Oro'
$ws.Range('B103').Value = 'Valle De Zaragoza'
$ws.Range('A122').Value = 'Ciudad De México'
$ws.Range('B136').Value = 'Coneto De Comonfort'
$ws.Range('B150').Value = 'Nombre De Dios'
$ws.Range('B153').Value = 'Pánuco De Coronado'
$ws.Range('B160').Value = 'San Juan Del Río'
$ws.Range('B161').Value = 'San Luis Del Cordero'
$ws.Range('A169').Value = 'Estado De México'
$ws.Range('B169').Value = 'Almoloya De Juárez'
$ws.Range('B174').Value = 'Chapa De Mota'
$ws.Range('B175').Value = 'Coacalco De Berriozábal'
$ws.Range('B178').Value = 'Ecatepec De Morelos'
$ws.Range('B180').Value = 'Ixtapan De La Sal'
$ws.Range('B184').Value = 'Naucalpan De Juárez'
$ws.Range('B188').Value = 'San Felipe Del Progreso'
$ws.Range('B195').Value = 'Tlalnepantla De Baz'
$ws.Range('A199').Value = 'Guanajuato'
$ws.Range('B201').Value = 'Apaseo El Alto'
$ws.Range('B205').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B216').Value = 'San Francisco Del Rincón'
$ws.Range('B217').Value = 'San Luis De La Paz'
$ws.Range('B219').Value = 'Silao De La Victoria'
$ws.Range('B226').Value = 'Acapulco De Juárez'
$ws.Range('B227').Value = 'Atenango Del Río'
$ws.Range('B229').Value = 'Atoyac De Álvarez'
$ws.Range('B230').Value = 'Ayutla De Los Libres'
$ws.Range('B232').Value = 'Chilapa De Álvarez'
$ws.Range('B233').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B235').Value = 'Coyuca De Benítez'
$ws.Range('B237').Value = 'Cutzamala De Pinzón'
$ws.Range('B241').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B249').Value = 'Técpan De Galeana'
$ws.Range('B251').Value = 'Tixtla De Guerero'
$ws.Range('B252').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B263').Value = 'Pachuca De Soto'
$ws.Range('B264').Value = 'Progreso De Obregón'
$ws.Range('B266').Value = 'Tula De Allende'
$ws.Range('B267').Value = 'Tulancingo De Bravo'
$ws.Range('B270').Value = 'Ahualulco De Mercado'
$ws.Range('B274').Value = 'Autlán De Navarro'
$ws.Range('B281').Value = 'Encarnación De Díaz'
$ws.Range('B285').Value = 'Huejuquilla El Alto'
$ws.Range('B286').Value = 'Ixtlahuacán Del Río'
$ws.Range('B290').Value = 'Lagos De Moreno'
$ws.Range('B293').Value = 'Ojuelos De Jalisco'
$ws.Range('B294').Value = 'San Cristóbal De La Barranca'
$ws.Range('B296').Value = 'San Juan De Los Lagos'
$ws.Range('B297').Value = 'San Juanito De Escobedo'
$ws.Range('B299').Value = 'Tamazula De Gordiano'
$ws.Range('B302').Value = 'Teocuitatlán De Corona'
$ws.Range('B303').Value = 'Tepatitlán De Morelos'
$ws.Range('B305').Value = 'Tizapán El Alto'
$ws.Range('B309').Value = 'Yahualica De González Gallo'
$ws.Range('B317').Value = 'Cojumatlán De Régules'
$ws.Range('B344').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B360').Value = 'Puente De Ixtla'
$ws.Range('B378').Value = 'Mier Y Noriega'
$ws.Range('B380').Value = 'San Nicolás De Los Garza'
$ws.Range('B385').Value = 'Guevea De Humboldt'
$ws.Range('B386').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B387').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B388').Value = 'Ixtlán De Juárez'
$ws.Range('B389').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B392').Value = 'Mariscala De Juárez'
$ws.Range('B393').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B394').Value = 'Oaxaca De Juárez'
$ws.Range('B410').Value = 'Teotitlán Del Valle'
$ws.Range('B411').Value = 'Tepelmeme Villa De Morelos'
$ws.Range('B412').Value = 'Villa Sola De Vega'
$ws.Range('B426').Value = 'Los Reyes De Juárez'
$ws.Range('B427').Value = 'Palmar De Bravo'
$ws.Range('B432').Value = 'San Salvador El Seco'
$ws.Range('B436').Value = 'Tepexi De Rodríguez'
$ws.Range('B437').Value = 'Tetela De Ocampo'
$ws.Range('B438').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B445').Value = 'Amealco De Bonfil'
$ws.Range('B446').Value = 'Cadereyta De Montes'
$ws.Range('B447').Value = 'Jalpan De Serra'
$ws.Range('B448').Value = 'Landa De Matamoros'
$ws.Range('B459').Value = 'Ciudad Del Maíz'
$ws.Range('B464').Value = 'San Ciro De Acosta'
$ws.Range('B466').Value = 'Santa María Del Río'
$ws.Range('B468').Value = 'Villa De Ramos'
$ws.Range('B500').Value = 'Nacozari De García'
$ws.Range('B522').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B533').Value = 'Cosamaloapan De Carpio'
$ws.Range('B537').Value = 'Ignacio De La Llave'
$ws.Range('B539').Value = 'Martínez De La Torre'
$ws.Range('B546').Value = 'Sayula De Alemán'
$ws.Range('B563').Value = 'Concepción Del Oro'
$ws.Range('B574').Value = 'Jiménez Del Teul'
$ws.Range('B579').Value = 'Nochistlán De Mejía'
$ws.Range('B585').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B587').Value = 'Villa De Cos'

# Fix floating point precision artifact
$ws.Range('D168').Value = 0.0929358044587698

# Remove footer note rows (595-599)
$ws.Range("A595:D599").EntireRow.Delete() | Out-Null

